$wb = $excel.ActiveWorkbook

# --- Norway ---
# Duplicate the "Turkey" sheet (it carries the formatting template used for
# the newer country sheets: no explicit row heights, wide column D, and the
# "select whole rows" selection state) and drop it at the end of the tab strip.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-2931/T3059"
$norway.Range("B2").Value = "Norway Market"

# --- Poland ---
$turkey2 = $wb.Worksheets.Item("Turkey")
$turkey2.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/3102"
$poland.Range("B2").Value = "Poland Market"

# The Norway tab is the one left active/selected after the edits.
$norway.Activate()
